# Apply the tracked "想去人数" (want-to-go count) refresh to 北京-漫展信息.xlsx
# Sheets: 1=展览 (Exhibition), 2=演出 (Performance), 3=本地生活 (Local life), 4=全部类型 (All types)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) 展览 sheet: straightforward F-column ("想去人数") value bumps.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @(
    @{ Row = 5;  Value = 5205 },
    @{ Row = 6;  Value = 5205 },
    @{ Row = 7;  Value = 142 },
    @{ Row = 9;  Value = 522 },
    @{ Row = 11; Value = 1172 },
    @{ Row = 12; Value = 732 },
    @{ Row = 13; Value = 5113 },
    @{ Row = 14; Value = 27 },
    @{ Row = 16; Value = 87 },
    @{ Row = 17; Value = 255 },
    @{ Row = 18; Value = 255 },
    @{ Row = 19; Value = 243 },
    @{ Row = 22; Value = 3866 },
    @{ Row = 23; Value = 43 },
    @{ Row = 24; Value = 3768 },
    @{ Row = 25; Value = 182 },
    @{ Row = 26; Value = 177 },
    @{ Row = 28; Value = 228 },
    @{ Row = 30; Value = 207 },
    @{ Row = 36; Value = 15 },
    @{ Row = 37; Value = 6698 },
    @{ Row = 38; Value = 1080 },
    @{ Row = 39; Value = 508 },
    @{ Row = 41; Value = 972 },
    @{ Row = 43; Value = 1366 },
    @{ Row = 44; Value = 171 },
    @{ Row = 45; Value = 677 },
    @{ Row = 47; Value = 2287 },
    @{ Row = 50; Value = 774 },
    @{ Row = 51; Value = 922 }
)
foreach ($u in $expoUpdates) {
    $wsExpo.Cells.Item($u.Row, 6).Value = $u.Value
}

# ---------------------------------------------------------------------------
# 2) 演出 sheet: a brand new performance (2024-07-21, 王子健2024巡回演出) is
#    inserted as row 16, pushing the former rows 16-26 down to 17-27. The
#    leading index column (A, "序号") keeps walking 0,1,2,... by row position
#    rather than traveling with the shifted content, so it is rewritten
#    afterwards back to row-number minus one for every data row.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Rows.Item(16).Insert()

# Copy formatting (bold + border style) from the neighbouring index cell so
# the new A16 matches the rest of the column.
$wsShow.Range("A17").Copy()
$wsShow.Range("A16").PasteSpecial(-4122)
$wsShow.Application.CutCopyMode = $false

$wsShow.Range("B16").Value = "2024-07-21"
$wsShow.Range("C16").Value = "北京·系统任务：重生之我是音乐一体机！王子健2024巡回演出"
$wsShow.Range("D16").Value = "北新桥街道板桥南巷7号人民美术印刷厂内 北京乐空间"
$wsShow.Range("E16").Value = "2024.07.21 20:00-07.21 22:00"
$wsShow.Range("F16").Value = 0
$wsShow.Range("G16").Value = 328
$wsShow.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=87587"
$wsShow.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202406/pW1onh2Z1718675834530.jpeg"

# Restore the sequential index column (A = row - 1) for every data row,
# since the insert above shifted it along with the rest of the row.
for ($r = 2; $r -le 27; $r++) {
    $wsShow.Cells.Item($r, 1).Value = $r - 1
}

# One of the shifted rows (now row 25, "Marcin Patrzalek ... 指弹吉他音乐会")
# also had its want-to-go count bumped from 810 to 811 in the same refresh.
$wsShow.Range("F25").Value = 811

# ---------------------------------------------------------------------------
# 3) 全部类型 sheet: same kind of F-column value bumps as 展览 (different row
#    offsets because this sheet merges every category sorted by date).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @(
    @{ Row = 7;  Value = 5205 },
    @{ Row = 8;  Value = 5205 },
    @{ Row = 9;  Value = 142 },
    @{ Row = 12; Value = 522 },
    @{ Row = 13; Value = 1172 },
    @{ Row = 14; Value = 732 },
    @{ Row = 15; Value = 5113 },
    @{ Row = 16; Value = 27 },
    @{ Row = 18; Value = 87 },
    @{ Row = 19; Value = 255 },
    @{ Row = 20; Value = 255 },
    @{ Row = 21; Value = 243 },
    @{ Row = 24; Value = 3866 },
    @{ Row = 25; Value = 3768 },
    @{ Row = 26; Value = 182 },
    @{ Row = 27; Value = 177 },
    @{ Row = 28; Value = 228 },
    @{ Row = 30; Value = 207 },
    @{ Row = 35; Value = 15 },
    @{ Row = 37; Value = 6698 },
    @{ Row = 38; Value = 1080 },
    @{ Row = 39; Value = 508 },
    @{ Row = 42; Value = 972 },
    @{ Row = 44; Value = 1367 },
    @{ Row = 45; Value = 171 },
    @{ Row = 46; Value = 677 },
    @{ Row = 47; Value = 2287 },
    @{ Row = 49; Value = 774 },
    @{ Row = 50; Value = 922 }
)
foreach ($u in $allUpdates) {
    $wsAll.Cells.Item($u.Row, 6).Value = $u.Value
}
